$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.889.25'
$ws.Range('E2').Value = '  -0.08%  '

$ws.Range('D3').Value = '3.120.00'
$ws.Range('E3').Value = '  +1.11%  '

$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.35'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.45%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '172.11'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.13%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.05%  '

$ws.Range('E9').Value = '  -2.53%  '

$ws.Range('E10').Value = '  -0.52%  '

$ws.Range('E11').Value = '  +0.10%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000247'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.99%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '37.08'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.55%  '

$ws.Range('E14').Value = '  -1.16%  '

$ws.Range('D15').Value = '3.634.24'
$ws.Range('E15').Value = '  +1.02%  '

$ws.Range('D16').Value = '66.860.13'
$ws.Range('E16').Value = '  -0.13%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.17'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.27%  '

$ws.Range('D18').Value = '3.114.74'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.28'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.62%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '478.48'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.84%  '

$ws.Range('B21').Value = 'Polygon'
$ws.Range('C21').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.713'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.47%  '

$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.97'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.83%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '84.07'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.93%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.35'
$ws.Range('D24').Style = 'Normal'

$ws.Range('E25').Value = '  -3.67%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.08'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.54%  '

$ws.Range('E27').Value = '  +0.00%  '

$ws.Range('E28').Value = '  -1.10%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.38'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.99%  '

$ws.Range('E30').Value = '  +0.52%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '28.60'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.61%  '

$ws.Range('E32').Value = '  +1.03%  '

$ws.Range('D33').Value = '0.0₃0943'
$ws.Range('E33').Value = '  -7.47%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.07%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.86'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.32%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.975'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.89%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '47.15'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.62%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '50.22'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.15%  '

$ws.Range('E39').Value = '  -3.08%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.311'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.45%  '

$ws.Range('E42').Value = '  -0.09%  '

$ws.Range('D43').Value = '2.823.93'
$ws.Range('E43').Value = '  +2.21%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '384.53'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.19%  '

$ws.Range('E45').Value = '  -1.54%  '

$ws.Range('E46').Value = '  -9.18%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '135.14'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.40%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '24.79'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.57%  '

$ws.Range('E50').Value = '  -1.93%  '

$ws.Range('E51').Value = '  -0.58%  '
